# Insert a new weekly price record as row 241 in the "Cilantro" price
# history table, pushing the existing rows 241-263 down to 242-264
# (dimension grows from A1:R263 to A1:R264).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 241..263 down by one, creating a blank row at 241.
$ws.Rows.Item(241).Insert()

# Populate the new row 241 with the new weekly observation.
$ws.Range("A241").Value = 8
$ws.Range("B241").Value = "Terminal La Palmera de La Serena"
$ws.Range("C241").Value = "Coquimbo"
$ws.Range("D241").Value = 45223
$ws.Range("E241").Value = 4
$ws.Range("F241").Value = 100112040
$ws.Range("G241").Value = "Cilantro"
$ws.Range("H241").Value = "Sin especificar"
$ws.Range("I241").Value = "Primera"
$ws.Range("J241").Value = 2000
$ws.Range("K241").Value = 1400
$ws.Range("L241").Value = 1500
$ws.Range("M241").Value = 1450
$ws.Range("N241").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O241").Value = "Provincia del Elquí"
$ws.Range("P241").Value = 967
$ws.Range("Q241").Value = 1.5
$ws.Range("R241").Value = "Hortaliza"
